$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "A11" (ClassID / NFTID) - fill in the real ibc class / NFT name
# for the "beauty004" evidence row.
# ----------------------------------------------------------------------
$wsA11 = $wb.Worksheets.Item("A11")
$wsA11.Range("A2").Value = "ibc/20E8BC52ADBA0BC5DF70645C1D0B99F06B9B6DD681C86EC3C578B9F582F0A680"
$wsA11.Range("B2").Value = "beauty004"

# ----------------------------------------------------------------------
# Sheet "A12" (ClassID / NFTID) - fill in the real ibc class / NFT name
# for the "beauty005" evidence row.
# ----------------------------------------------------------------------
$wsA12 = $wb.Worksheets.Item("A12")
$wsA12.Range("A2").Value = "ibc/2DFDE04A57F9C0E738BB0E33DD65C62DE7D629E31C64F31975D864AAF457F4FF"
$wsA12.Range("B2").Value = "beauty005"

# ----------------------------------------------------------------------
# Sheet "A19" (TxHash / ChainID) - add the real tx hashes / chain ids,
# growing the table from 4 rows to 7 rows.
# ----------------------------------------------------------------------
$wsA19 = $wb.Worksheets.Item("A19")

$wsA19.Range("A2").Value = "860EF65C2EE998752A692D6D4A856F1BAD45275C76ED929FB445AC12D21776EB"
$wsA19.Range("A3").Value = "FE64BCBA489CC6E4830E112CFECF0863C0CBD9F76D6FF6019BA91F80F34046EA"
$wsA19.Range("A4").Value = "A8628074064986926D13484C3DBCFB24A5BE9C27A75E3EFD6770F10E0D6EA09D"

$wsA19.Range("B2").Value = "gon-irishub-1"
$wsA19.Range("B3").Value = "elgafar-1"
$wsA19.Range("B4").Value = "uni-6"
$wsA19.Range("B5").Value = "uptick_7000-2"
$wsA19.Range("B6").Value = "uni-6"
$wsA19.Range("B7").Value = "elgafar-1"

# New rows 5-7 don't inherit the sheet's row height (14pt, like rows 1-4)
# automatically, so set it explicitly to keep the table visually uniform.
$wsA19.Rows.Item(5).RowHeight = 14
$wsA19.Rows.Item(6).RowHeight = 14
$wsA19.Rows.Item(7).RowHeight = 14

# Match the font of the newly-created rows' B cells to the rest of the
# column before left-aligning them, so they collapse onto the same style
# record instead of spawning one new xf per distinct previous style.
$wsA19.Range("B3:B7").Font.Name = $wsA19.Range("B2").Font.Name
$wsA19.Range("B3:B7").Font.Size = $wsA19.Range("B2").Font.Size
$wsA19.Range("B3:B7").Font.Color = $wsA19.Range("B2").Font.Color

$wsA19.Range("B2").HorizontalAlignment = -4131
$wsA19.Range("B3").HorizontalAlignment = -4131
$wsA19.Range("B4").HorizontalAlignment = -4131
$wsA19.Range("B5").HorizontalAlignment = -4131
$wsA19.Range("B6").HorizontalAlignment = -4131
$wsA19.Range("B7").HorizontalAlignment = -4131

# ----------------------------------------------------------------------
# Selections / active sheet: A11 loses the tab selection, A12 and A19
# pick up a remembered cursor position, and A19 becomes the active tab
# (it is selected last).
# ----------------------------------------------------------------------
[void]$wsA11.Range("B6").Select()
[void]$wsA12.Range("C4").Select()
[void]$wsA19.Range("C8").Select()
